# "9th Stab - Cosmetic Changes"
#
# The watch-sheet keeps a rolling window of "current" rating columns.
# Every time the sheet is refreshed, a brand-new "current" column is
# inserted right after the firm/status columns (A, B) and the previous
# "current" column(s) slide right to become history. The newest column
# starts out unfilled ("UN") for every firm until the next real pull
# populates it; only the header row gets the new pull date right away.
#
# This pass rolls the window forward twice - once for Jun_15 and once
# for Jun_17 - which is why two fresh columns (C and D) appear and the
# formerly-current Jun_13 column (plus its already-collected ratings)
# shifts two slots to the right, landing in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $firstRow + $usedRange.Rows.Count - 1
$headerRow = $firstRow

$newDates = @("Jun_15", "Jun_17")

foreach ($newDate in $newDates) {
    # Insert a fresh "current" column right before the existing one
    # (column B), pushing every older column one slot to the right.
    $ws.Columns("B:B").Insert()

    # The header row records which pull this column represents...
    $ws.Cells.Item($headerRow, 2).Value = $newDate

    # ...while every firm's row starts this column as "UN" (unrated)
    # until the data for that date actually comes in.
    for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 2).Value = "UN"
    }
}

# Keep the now three "current-ish" columns (C, D, E) all the same
# cosmetic width as the original column (8 characters).
$ws.Columns("C:E").ColumnWidth = 7.1666666666667
